$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Etherium Cell', ['Token Artifact', '{T}, Sacrifice this artifact: Add one mana of any color.'])"
$ws.Range("A3").Value = "('Gremlin', ['Token Creature — Gremlin', '2/2'])"
$ws.Range("A4").Value = "('Ragavan', ['Token Legendary Creature — Monkey', '2/1'])"
$ws.Range("A5").Value = "('Tezzeret the Schemer Emblem', ['Emblem — Tezzeret', 'At the beginning of combat on your turn, target artifact you control becomes an artifact creature with base power and toughness 5/5.'])"

$ws.Rows("6:13").Delete()
